# CWR_Checklist_Template edits:
#  - CK_Identification: convert J9 to a text-formatted date string, add rows 10-12
#    (new test records, Allium / TEST-2 data) with hyperlinks on C10:C12 / I12
#  - CK_Crossability: add row 10 (more crossability data)
#  - CK_Threats: add rows 10-11 (Torrimpietra / Fregene threat records)
#  - Active sheet / selection moves from CK_Crossability to CK_Threats

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CK_Identification (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CK_Identification")

# J9 used to hold a date serial (42069, formatted m/d/yyyy). It becomes a
# free-text date string instead, so switch the cell to Text format first.
$ws2.Range("J9").NumberFormat = "@"
$ws2.Range("J9").Value = "20150308"

# --- Row 10 --------------------------------------------------------------
$ws2.Range("B10").Value = "Test dataset"
$ws2.Range("C10").Value = "http://pgrdiversity.bioversityinternational.org"
$ws2.Range("D10").Value = "This is a trial dataset used to test template import."
$ws2.Range("E10").Value = "TEST-1"
$ws2.Range("J10").NumberFormat = "@"
$ws2.Range("K10").Value = "IT-RM"
$ws2.Range("L10").Value = 201504
$ws2.Range("N10").Value = 3
$ws2.Range("T10").Value = "Allium"
$ws2.Range("U10").Value = "aestivum"
$ws2.Range("Y10").Value = "Allium aestivum"
$ws2.Range("AS10").Value = "Missing index reference"

# --- Row 11 ----------------------------------------------------------------
$ws2.Range("A11").Value = 3
$ws2.Range("B11").Value = "Test dataset"
$ws2.Range("C11").Value = "http://pgrdiversity.bioversityinternational.org"
$ws2.Range("D11").Value = "This is a trial dataset used to test template import."
$ws2.Range("E11").Value = "TEST-1"
$ws2.Range("J11").NumberFormat = "@"
$ws2.Range("K11").Value = "IT-RM"
$ws2.Range("L11").Value = " "
$ws2.Range("M11").Value = "ITA406"

# --- Row 12 ------------------------------------------------------------------
$ws2.Range("A12").Value = 4
$ws2.Range("B12").Value = "Test dataset"
$ws2.Range("C12").Value = "http://pgrdiversity.bioversityinternational.org"
$ws2.Range("D12").Value = "This is a trial dataset used to test template import."
$ws2.Range("E12").Value = "TEST-2"
$ws2.Range("F12").Value = "IT-RM"
$ws2.Range("G12").Value = "ITA406"
$ws2.Range("I12").Value = "http://bioversityinternational.org"
$ws2.Range("J12").NumberFormat = "@"
$ws2.Range("J12").Value = "20150310"
$ws2.Range("K12").Value = "IT-RM"
$ws2.Range("L12").Value = 201504
$ws2.Range("M12").Value = "ITA406"
$ws2.Range("N12").Value = 1
$ws2.Range("O12").Value = "Plantae"
$ws2.Range("T12").Value = "Triticum"
$ws2.Range("U12").Value = "aestivum"
$ws2.Range("V12").Value = "L."
$ws2.Range("W12").Value = "var. pippolense"
$ws2.Range("Y12").Value = "Triticum aestivum L. var. pippolense"
$ws2.Range("AD12").Value = "en@phoney wheat;it@grano finto"
$ws2.Range("AO12").Value = "Wheat"
$ws2.Range("AP12").Value = 20
$ws2.Range("AS12").Value = "Test 2 record"

# Hyperlinks for the newly added rows (same external targets as the row 9
# examples already on this sheet)
$ws2.Hyperlinks.Add($ws2.Range("C10"), "http://pgrdiversity.bioversityinternational.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C11"), "http://pgrdiversity.bioversityinternational.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C12"), "http://pgrdiversity.bioversityinternational.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I12"), "http://bioversityinternational.org") | Out-Null

# ---------------------------------------------------------------------------
# CK_Crossability (sheet3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CK_Crossability")

$ws3.Range("A10").Value = 4
$ws3.Range("B10").Value = "Triticum vulgare;Triticum muticum"
$ws3.Range("C10").Value = "forced"
$ws3.Range("E10").Value = 60

# ---------------------------------------------------------------------------
# CK_Threats (sheet4)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("CK_Threats")

# --- Row 10 ------------------------------------------------------------------
$ws4.Range("A10").Value = 3
$ws4.Range("B10").Value = 2
$ws4.Range("C10").Value = "Torrimpietra"
$ws4.Range("D10").Value = "DD"
$ws4.Range("E10").Value = "C"

# --- Row 11 ------------------------------------------------------------------
$ws4.Range("A11").Value = 4
$ws4.Range("B11").Value = 2
$ws4.Range("C11").Value = "Fregene"
$ws4.Range("D11").Value = "NE"
$ws4.Range("E11").Value = "D"
$ws4.Range("G11").Value = "FREG-1"
$ws4.Range("H11").Value = "2.2.2"
$ws4.Range("I11").Value = 3
$ws4.Range("J11").Value = 2015
$ws4.Range("M11").Value = "G4"
$ws4.Range("S11").Value = 410
$ws4.Range("T11").Value = "Some notes"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping.
# Each sheet's cursor moves one row below its newly added data, and the
# workbook's active tab moves from CK_Crossability to CK_Threats.
# ---------------------------------------------------------------------------
$ws2.Range("A13").Select()
$ws3.Range("A11").Select()
$ws4.Range("A12").Select()
$ws4.Activate()

Write-Host "edits applied"
